$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.105.02"
Set-TextValue "E2" "  -0.42%  "
Set-TextValue "D3" "1.627.49"
Set-TextValue "E3" "  -1.19%  "
Set-TextValue "E4" "  +0.07%  "
Set-TextValue "D5" "216.20"
Set-TextValue "E5" "  -1.14%  "
Set-TextValue "E6" "  +0.59%  "
Set-TextValue "E7" "  +0.10%  "
Set-TextValue "E8" "  -1.63%  "
Set-TextValue "D9" "0.0622"
Set-TextValue "D10" "20.04"
Set-TextValue "E10" "  -0.10%  "
Set-TextValue "D11" "0.0851"
Set-TextValue "E11" "  +0.51%  "
Set-TextValue "D12" "1.632.94"
Set-TextValue "E12" "  -1.04%  "
Set-TextValue "E13" "  -0.92%  "
Set-TextValue "E14" "  +0.00%  "
Set-TextValue "D15" "27.089.80"
Set-TextValue "E15" "  -0.44%  "
Set-TextValue "D16" "64.62"
Set-TextValue "E16" "  -4.31%  "
Set-TextValue "E17" "  -1.50%  "
Set-TextValue "D18" "213.48"
Set-TextValue "E18" "  -2.85%  "
Set-TextValue "E19" "  +0.04%  "
Set-TextValue "D20" "6.86"
Set-TextValue "E20" "  +0.84%  "
Set-TextValue "D21" "4.38"
Set-TextValue "E21" "  -1.40%  "
Set-TextValue "D22" "2.49"
Set-TextValue "E22" "  +0.82%  "
Set-TextValue "D23" "9.03"
Set-TextValue "E23" "  -2.04%  "
Set-TextValue "D24" "148.12"
Set-TextValue "E24" "  -0.01%  "
Set-TextValue "E25" "  +0.21%  "
Set-TextValue "D26" "7.28"
Set-TextValue "E26" "  -1.89%  "
Set-TextValue "E27" "  -1.10%  "
Set-TextValue "D28" "15.50"
Set-TextValue "E28" "  -1.96%  "
Set-TextValue "D29" "0.0503"
Set-TextValue "E29" "  -0.66%  "
Set-TextValue "E30" "  -0.97%  "
Set-TextValue "D31" "3.38"
Set-TextValue "E31" "  +0.25%  "
Set-TextValue "E32" "  -1.50%  "
Set-TextValue "D33" "1.311.17"
Set-TextValue "E33" "  +3.44%  "
Set-TextValue "D34" "1.55"
Set-TextValue "E34" "  -2.25%  "
Set-TextValue "E35" "  -0.11%  "
Set-TextValue "D36" "0.0174"
Set-TextValue "E36" "  -2.10%  "
Set-TextValue "D37" "0.840"
Set-TextValue "E37" "  -0.68%  "
Set-TextValue "E38" "  -1.57%  "
Set-TextValue "E39" "  +0.06%  "
Set-TextValue "D40" "2.25"
Set-TextValue "E40" "  +0.80%  "
Set-TextValue "E41" "  -1.00%  "
Set-TextValue "D42" "63.41"
Set-TextValue "E42" "  +1.66%  "
Set-TextValue "B43" "RocketPoolETH"
Set-TextValue "C43" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D43" "1.765.08"
Set-TextValue "E43" "  -1.31%  "
Set-TextValue "B44" "FraxShare"
Set-TextValue "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D44" "5.24"
Set-TextValue "E44" "  -2.74%  "
Set-TextValue "D45" "90.43"
Set-TextValue "E45" "  -1.79%  "
Set-TextValue "E46" "  -0.62%  "
Set-TextValue "E47" "  +17.43%  "
Set-TextValue "D48" "0.792"
Set-TextValue "E48" "  +16.81%  "
Set-TextValue "D49" "0.0514"
Set-TextValue "E49" "  +0.16%  "
Set-TextValue "D50" "7.51"
Set-TextValue "E50" "  -2.64%  "
Set-TextValue "E51" "  +0.09%  "
